# Insert a new data row at row 337 (pushing the existing rows 337:420 down
# to 338:421) on the single worksheet of the "Hortaliza, Macroferia Regional
# de Talca - Repollo" workbook, and populate the new row with the new
# weekly price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 337:420 down to 338:421, leaving a blank row 337 (formatting
# of column D's date style is carried down automatically by Excel).
$ws.Rows.Item(337).Insert()

# Fill the newly-inserted row 337 with the new record.
$ws.Range("A337").Value = 5
$ws.Range("B337").Value = "Macroferia Regional de Talca"
$ws.Range("C337").Value = "Maule"
$ws.Range("D337").Value = 44932
$ws.Range("E337").Value = 7
$ws.Range("F337").Value = 100112006
$ws.Range("G337").Value = "Repollo"
$ws.Range("H337").Value = "Crespo record"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 3000
$ws.Range("K337").Value = 1000
$ws.Range("L337").Value = 1000
$ws.Range("M337").Value = 1000
$ws.Range("N337").Value = "$/unidad"
$ws.Range("O337").Value = "Región del Maule"
$ws.Range("P337").Value = 1000
$ws.Range("Q337").Value = 1
$ws.Range("R337").Value = "Hortaliza"
